$wb = $excel.ActiveWorkbook

function Set-RowData {
    param(
        $ws,
        [int]$row,
        [int]$a,
        [string]$b,
        [string]$c,
        [string]$d,
        [string]$e,
        [double]$f,
        [double]$g,
        [string]$h,
        [string]$i
    )
    $ws.Cells.Item($row, 1).Value = $a
    # "2024-05-25" reads as an ISO date, so a bare assignment gets
    # auto-coerced into an Excel date serial; a leading apostrophe forces
    # it to stay literal text (matching the sheet's existing B-column
    # cells, which are all stored as plain text). Re-apply the "Normal"
    # style afterwards so the quote-prefix flag doesn't linger as a
    # visible format difference on the cell.
    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value = "'" + $b
    $bCell.Style = "Normal"
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
}

function Format-NewRowHeaderCell {
    param($ws, [int]$row)
    $c = $ws.Cells.Item($row, 1)
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
    $c.Borders.LineStyle = 1
}

# ---- Sheet "展览" (index 1) ----
$ws1 = $wb.Worksheets.Item(1)

# A new entry (2024-05-25, 南宁·第五人格Only1.0) is inserted as row 5,
# pushing the existing rows 5-10 down to 6-11.
$ws1.Rows.Item(5).Insert()
Format-NewRowHeaderCell $ws1 5
Set-RowData $ws1 5 4 "2024-05-25" "南宁·第五人格Only1.0" "新阳路227号南宁第三人民医院旁新秀佳园对面 卡尔顿东方银龙酒店" "2024.05.25 10:00-05.25 17:30" 4 68 "https://show.bilibili.com/platform/detail.html?id=84954" "//i0.hdslb.com/bfs/openplatform/202404/w5iZT4wE1714189905443.jpeg"

# Renumber the serial-number column (A) for the rows that got shifted down.
for ($r = 6; $r -le 11; $r++) {
    $ws1.Cells.Item($r, 1).Value = $r - 1
}

# Small "want-to-go" count refreshes on pre-existing rows.
$ws1.Cells.Item(2, 6).Value = 7533
$ws1.Cells.Item(4, 6).Value = 213
$ws1.Cells.Item(6, 6).Value = 248
$ws1.Cells.Item(7, 6).Value = 1134
$ws1.Cells.Item(10, 6).Value = 135
$ws1.Cells.Item(11, 6).Value = 35

# ---- Sheet "全部类型" (index 4) ----
$ws4 = $wb.Worksheets.Item(4)

$ws4.Rows.Item(5).Insert()
Format-NewRowHeaderCell $ws4 5
Set-RowData $ws4 5 4 "2024-05-25" "南宁·第五人格Only1.0" "新阳路227号南宁第三人民医院旁新秀佳园对面 卡尔顿东方银龙酒店" "2024.05.25 10:00-05.25 17:30" 4 68 "https://show.bilibili.com/platform/detail.html?id=84954" "//i0.hdslb.com/bfs/openplatform/202404/w5iZT4wE1714189905443.jpeg"

for ($r = 6; $r -le 12; $r++) {
    $ws4.Cells.Item($r, 1).Value = $r - 1
}

$ws4.Cells.Item(2, 6).Value = 7533
$ws4.Cells.Item(4, 6).Value = 213
$ws4.Cells.Item(6, 6).Value = 248
$ws4.Cells.Item(7, 6).Value = 1134
$ws4.Cells.Item(11, 6).Value = 135
$ws4.Cells.Item(12, 6).Value = 35

Write-Output "done"
